$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are plain numeric-looking strings in the source data
# (e.g. "42.739.83", "309.28"). Excel auto-converts plain numeric strings to
# real numbers on assignment, which would change the stored cell type from the
# original text. Force text storage by temporarily marking the cell as Text,
# then clearing the format flag again so the cell keeps its original (default)
# style while the stored value remains a text string.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "42.739.83"
$ws.Range("E2").Value = "  -0.88%  "
Set-TextValue $ws.Range("D3") "2.528.71"
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "309.28"
$ws.Range("E5").Value = "  -1.87%  "
Set-TextValue $ws.Range("D6") "101.89"
$ws.Range("E6").Value = "  +4.37%  "
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E9").Value = "  -1.47%  "
Set-TextValue $ws.Range("D10") "35.98"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  -0.93%  "
Set-TextValue $ws.Range("D12") "7.33"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("E13").Value = "  +0.25%  "
Set-TextValue $ws.Range("D14") "2.918.30"
$ws.Range("E14").Value = "  -2.27%  "
Set-TextValue $ws.Range("D15") "15.73"
$ws.Range("E15").Value = "  +3.02%  "
Set-TextValue $ws.Range("D16") "2.515.15"
$ws.Range("E16").Value = "  +0.32%  "
Set-TextValue $ws.Range("D17") "0.810"
$ws.Range("E17").Value = "  -4.21%  "
Set-TextValue $ws.Range("D18") "42.705.77"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  -4.50%  "
Set-TextValue $ws.Range("D22") "69.39"
$ws.Range("E22").Value = "  -0.11%  "
Set-TextValue $ws.Range("D23") "244.29"
$ws.Range("E23").Value = "  -2.73%  "
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("E26").Value = "  +0.04%  "
Set-TextValue $ws.Range("D27") "26.17"
$ws.Range("E27").Value = "  -4.24%  "
$ws.Range("E28").Value = "  -3.66%  "
Set-TextValue $ws.Range("D29") "39.14"
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("E30").Value = "  -1.33%  "
Set-TextValue $ws.Range("D31") "157.36"
$ws.Range("E31").Value = "  +0.30%  "
Set-TextValue $ws.Range("D32") "5.79"
$ws.Range("E32").Value = "  -1.01%  "
Set-TextValue $ws.Range("D33") "2.81"
$ws.Range("E33").Value = "  +11.88%  "
Set-TextValue $ws.Range("D34") "0.0790"
$ws.Range("E34").Value = "  -1.91%  "
Set-TextValue $ws.Range("D35") "2.63"
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("E36").Value = "  -5.35%  "
$ws.Range("E37").Value = "  -6.84%  "
Set-TextValue $ws.Range("D38") "18.13"
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("E41").Value = "  +5.73%  "
Set-TextValue $ws.Range("D42") "22.07"
$ws.Range("E42").Value = "  -5.41%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("E45").Value = "  -1.55%  "
Set-TextValue $ws.Range("D46") "1.994.07"
$ws.Range("E46").Value = "  -0.89%  "
Set-TextValue $ws.Range("D47") "8.89"
$ws.Range("E47").Value = "  -1.13%  "
Set-TextValue $ws.Range("D48") "2.768.93"
$ws.Range("E48").Value = "  -2.45%  "
Set-TextValue $ws.Range("D49") "80.40"
$ws.Range("E49").Value = "  -3.53%  "
$ws.Range("E50").Value = "  -2.58%  "
Set-TextValue $ws.Range("D51") "72.17"
$ws.Range("E51").Value = "  -3.00%  "
